# Add a new closing slide ("Благодаря за вниманието" / "Thank you for your
# attention") at the end of the deck, using the same "Title and Content"
# layout as the preceding slide, but with only the title placeholder
# populated (no body content placeholder).
#
# We build it by duplicating the last existing slide (which already uses
# the right layout and carries the standard PowerPoint slide boilerplate -
# xfrm, clrMapOvr, creationId, ...), then stripping it down to just the
# title placeholder and updating its text.

$p = $ppt.ActivePresentation

$lastIndex = $p.Slides.Count
$sourceSlide = $p.Slides.Item($lastIndex)

# Duplicate the last slide; the duplicate is inserted immediately after it,
# i.e. becomes the new last slide. Duplicate() returns a SlideRange - grab
# the new Slide from it.
$dup = $sourceSlide.Duplicate()
$newSlide = $dup.Item(1)

# The duplicated slide carries over the source's hyperlink relationship
# (on its content placeholder). Clear it before removing the shape so the
# now-unused relationship doesn't linger in the new slide's part.
if ($newSlide.Hyperlinks.Count -gt 0) {
    for ($i = $newSlide.Hyperlinks.Count; $i -ge 1; $i--) {
        $newSlide.Hyperlinks.Item($i).Address = ""
    }
}

# Remove every shape except the title placeholder (shape 1), leaving just
# the title. Placeholder shapes need a couple of Delete() calls in this
# host before they are actually removed rather than just reset.
for ($i = $newSlide.Shapes.Count; $i -ge 2; $i--) {
    $newSlide.Shapes.Item($i).Delete()
    $newSlide.Shapes.Item($i).Delete()
}

# Set the title text.
$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Благодаря за вниманието"
